$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) cells for affected rows.
# D-column values are assigned with a leading apostrophe so Excel keeps
# them as literal text (matching the original inline-string cells) instead
# of auto-parsing them as numbers/dates.

$ws.Range("D2").Value = "'54.324.05"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "'2.285.56"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'495.65"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "'127.55"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "'2.282.30"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "'0.0941"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").Value = "'2.674.93"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "'21.68"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'54.145.66"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "'2.295.62"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'4.04"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'303.06"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "'6.31"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'63.67"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'2.392.52"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "'164.22"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").Value = "'1.60"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'0.0₃0680"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'5.87"
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "'1.07"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'17.52"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'0.871"
$ws.Range("E39").Value = "  +6.16%  "
$ws.Range("D40").Value = "'3.64"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'35.29"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'3.34"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D47").Value = "'0.0890"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'0.546"
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").Value = "'238.01"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'0.0204"
$ws.Range("E51").Value = "  -0.38%  "

# Rows 45 and 46: RenderToken and Aave swap list positions (the rank index
# stored in column A is unaffected by the swap).
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'4.92"
$ws.Range("E45").Value = "  +4.47%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'126.30"
$ws.Range("E46").Value = "  -0.34%  "
